$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new Price (Column D) value. Omitted rows keep existing D value.
$priceUpdates = @{
    2  = "44.016.94"
    3  = "2.237.48"
    4  = "1.01"
    5  = "304.44"
    6  = "96.18"
    7  = "0.571"
    9  = "0.522"
    10 = "34.45"
    11 = "0.0802"
    12 = "7.13"
    14 = "2.583.18"
    15 = "2.245.95"
    16 = "0.817"
    17 = "13.55"
    18 = "43.846.27"
    19 = "0.0₃0956"
    20 = "12.26"
    21 = "6.19"
    22 = "64.55"
    23 = "237.59"
    24 = "2.91"
    26 = "1.93"
    27 = "9.88"
    28 = "2.12"
    29 = "35.92"
    30 = "5.94"
    31 = "19.88"
    32 = "153.34"
    33 = "3.36"
    34 = "0.0801"
    36 = "0.118"
    38 = "1.75"
    39 = "14.91"
    40 = "3.77"
    41 = "0.0299"
    42 = "3.27"
    44 = "1.739.07"
    45 = "85.51"
    46 = "5.08"
    47 = "0.187"
    48 = "99.55"
    49 = "8.14"
    50 = "54.03"
    51 = "67.63"
}

# Map of row -> new Volume(1h) value (Column E).
$volumeUpdates = @{
    2  = "  +0.57%  "
    3  = "  -0.40%  "
    4  = "  +0.33%  "
    5  = "  -5.73%  "
    6  = "  -5.32%  "
    7  = "  -1.74%  "
    8  = "  +0.26%  "
    9  = "  -6.33%  "
    10 = "  -7.33%  "
    11 = "  -3.61%  "
    12 = "  -7.78%  "
    13 = "  -2.88%  "
    14 = "  -0.20%  "
    15 = "  -0.13%  "
    16 = "  -4.66%  "
    17 = "  -4.47%  "
    18 = "  +0.38%  "
    19 = "  -3.24%  "
    20 = "  -9.68%  "
    21 = "  -5.58%  "
    22 = "  -1.32%  "
    23 = "  +0.56%  "
    24 = "  -8.14%  "
    25 = "  +0.10%  "
    26 = "  -10.05%  "
    27 = "  -2.65%  "
    28 = "  -3.11%  "
    29 = "  -3.35%  "
    30 = "  -6.02%  "
    31 = "  -1.42%  "
    32 = "  -4.82%  "
    33 = "  +9.83%  "
    34 = "  -5.87%  "
    35 = "  -1.56%  "
    36 = "  -0.53%  "
    37 = "  -8.43%  "
    38 = "  -9.05%  "
    39 = "  -7.41%  "
    40 = "  -11.30%  "
    41 = "  -5.71%  "
    42 = "  -12.99%  "
    44 = "  -4.15%  "
    45 = "  +3.67%  "
    46 = "  -2.57%  "
    47 = "  -6.09%  "
    48 = "  -3.61%  "
    49 = "  -3.46%  "
    50 = "  -7.90%  "
    51 = "  -11.59%  "
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}
